$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Head Alpha Pro"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "179"
$ws.Range("A3").Value = "Head Alpha Control"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "129"

$ws.Range("B2").Select()
